$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '51.670.46'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.923.18'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.93%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '353.46'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.76'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.559'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.629'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.99'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.04%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0894'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +4.74%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.137'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.65'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.07%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.93'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.379.95'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.83%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.915.36'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.31%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.984'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '51.713.90'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.69%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.60'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.21'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +4.93%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '3.24'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -3.88%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0980'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '70.80'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '269.65'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.78'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.181'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +9.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '27.14'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +3.07%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.47'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +19.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.107'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +16.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '10.60'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.19'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '36.87'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -5.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '52.29'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.42%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0438'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -4.68%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.89'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -16.33%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.25'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -1.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.26'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -1.72%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.02'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.32%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.71'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.118'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.96%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '23.11'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +3.63%  '
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.16%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '114.58'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -5.86%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.45'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -4.29%  '
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.134.95'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.03%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0327'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.21%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '9.22'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +2.42%  '
